# Apply the changes described by the diff to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column D (emissions_intensity) values for the listed rows.
$ws.Range("D5").Value = 150
$ws.Range("D8").Value = 6
$ws.Range("D10").Value = 6
$ws.Range("D14").Value = 9
$ws.Range("D16").Value = 10
$ws.Range("D19").Value = 6
$ws.Range("D23").Value = 10
$ws.Range("D26").Value = 5
$ws.Range("D27").Value = 4
$ws.Range("D30").Value = 7
$ws.Range("D31").Value = 7
$ws.Range("D32").Value = 10

# Best-fit column A so it fits the longest entry ("Plant Power Fast Food"),
# matching the ~20.7-character-wide column from the diff.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(1).ColumnWidth = 19.833333333333332

# Update the active selection to I23 (no multi-cell selection anymore).
$ws.Range("I23").Select() | Out-Null
